## feat: add 2022-Q4 data
##
## The existing "2022-Q2" tab (2nd sheet, sheetId 2) keeps its sheetId/rId but is
## repointed to hold the new 2022-Q4 fund-holdings table; a fresh tab -- an exact
## copy of the original "2022-Q2" content -- is appended after it (sheetId 3) and
## renamed back to "2022-Q2". The summary tab ("总计") gets a new row 2 for
## 2022-Q4, pushing the existing 2022-Q2 summary row down to row 3.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# --- 1. Clone the current "2022-Q2" sheet (inserted right after itself) so its
#        data survives under a new trailing tab; $q2 keeps its original
#        sheetId/rId in place, the clone picks up a fresh one ---
$q2.Copy([System.Reflection.Missing]::Value, $q2)
$q2clone = $wb.Worksheets.Item(3)

# --- 2. Rename the original sheet away first (freeing up the "2022-Q2" name),
#        then rename the clone back to "2022-Q2" ---
$q2.Name = "2022-Q4"
$q2clone.Name = "2022-Q2"

# --- 3. Wipe the original sheet and refill it with the 2022-Q4 fund table ---
$q2.Cells.Clear()

$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"

# Header row + the A-column index cells reuse the bold/boxed style from the
# "总计" sheet (style index 2 in the original workbook).
$summary.Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q2.Range("A2:A7").PasteSpecial(-4122)

$q4rows = @(
    @(0, "580006", "东吴新经济混合A",         "1.16", "89.15", "4.46", "0.0517", 7),
    @(1, "012617", "东吴新经济混合C",         "0.42", "89.15", "4.46", "0.0187", 7),
    @(2, "009856", "中加新兴成长混合C",       "0.46", "93.78", "3.40", "0.0156", 4),
    @(3, "009855", "中加新兴成长混合A",       "0.19", "93.78", "3.40", "0.0065", 4),
    @(4, "007315", "汇安嘉盈一年持有期债券A", "0.16", "24.33", "1.18", "0.0019", 6),
    @(5, "010270", "汇安嘉盈一年持有期债券C", "0.15", "24.33", "1.18", "0.0018", 6)
)

# Fund code (B) and the numeric-looking text columns (D:G) must stay text --
# mark them "@" up front so the leading zeros in the fund codes (e.g. 012617)
# and the literal "1.16"/"0.0517" strings survive instead of being coerced to
# numbers.
$q2.Range("B2:B7").NumberFormat = "@"
$q2.Range("D2:G7").NumberFormat = "@"

for ($i = 0; $i -lt $q4rows.Count; $i++) {
    $r = 2 + $i
    $row = $q4rows[$i]
    $q2.Range("A$r").Value = $row[0]
    $q2.Range("B$r").Value = $row[1]
    $q2.Range("C$r").Value = $row[2]
    $q2.Range("D$r").Value = $row[3]
    $q2.Range("E$r").Value = $row[4]
    $q2.Range("F$r").Value = $row[5]
    $q2.Range("G$r").Value = $row[6]
    $q2.Range("H$r").Value = $row[7]
}

# --- 4. "总计": push the existing 2022-Q2 summary row down to row 3, and write
#        the new 2022-Q4 summary row into row 2 ---
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.07000000000000001

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.1

# --- 5. Leave the workbook's active tab on "总计", same as before the edit ---
$summary.Activate()

Write-Host "2022-Q4 sheet added"
